$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Core data change: "TC_G" (growth rate) for 2020 goes from 0.9% to 9%
$ws.Range("N18").Value = 0.09

# New running-total row under the 2020 column, formatted like the neighboring
# currency cells (row 18/22) instead of the previous blank "General" style
$ws.Range("B19").Formula = "=SUM(C18:H18)"
$ws.Range("B19").NumberFormat = $ws.Range("C18").NumberFormat

Write-Output "edit applied"
